$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Improvements" comment for row 9 (C06 / image orientation fix)
$ws.Range("D9").Value = "Mit der Image orientation verbesserung gibt es wenigere Fehler"

# Update the selected cell on the sheet view
$ws.Range("C12").Select()
